$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row - Right column
$ws.Range("B11").Value = 5

# Update "Total" row - Right column and the Correct/Total text in Max column
$ws.Range("B12").Value = 65
$ws.Range("E12").Value = "65/140"
